$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data range entirely (A1:E12) before rewriting
$ws.Range("A1:F12").ClearContents()

# Row 1 stays empty.

# Row 2: headers
$ws.Cells.Item(2,1).Value = "Player Name"
$ws.Cells.Item(2,2).Value = "Gender"
$ws.Cells.Item(2,3).Value = "Final Score"
$ws.Cells.Item(2,4).Value = "Danger Level"
$ws.Cells.Item(2,5).Value = "Game Outcome"
$ws.Cells.Item(2,6).Value = "Timestamp"

# Row 3
$ws.Cells.Item(3,1).Value = "Siam"
$ws.Cells.Item(3,2).Value = "Male"
$ws.Cells.Item(3,3).Value = 50
$ws.Cells.Item(3,4).Value = 1
$ws.Cells.Item(3,5).Value = "Game Over"
$ws.Cells.Item(3,6).Value = "2025-05-24 14:54:39"

# Row 4
$ws.Cells.Item(4,1).Value = "Abir"
$ws.Cells.Item(4,2).Value = "Male"
$ws.Cells.Item(4,3).Value = 80
$ws.Cells.Item(4,4).Value = 1
$ws.Cells.Item(4,5).Value = "Quit"
$ws.Cells.Item(4,6).Value = "2025-05-24 14:57:13"

# Row 5
$ws.Cells.Item(5,1).Value = "Hujaifa"
$ws.Cells.Item(5,2).Value = "Male"
$ws.Cells.Item(5,3).Value = 80
$ws.Cells.Item(5,4).Value = 1
$ws.Cells.Item(5,5).Value = "Game Over"
$ws.Cells.Item(5,6).Value = "2025-05-24 14:58:25"

# Row 6
$ws.Cells.Item(6,1).Value = "Pial"
$ws.Cells.Item(6,2).Value = "Female"
$ws.Cells.Item(6,3).Value = 410
$ws.Cells.Item(6,4).Value = 3
$ws.Cells.Item(6,5).Value = "Game Over"
$ws.Cells.Item(6,6).Value = "2025-05-24 14:59:35"

# Column widths
$ws.Columns.Item(5).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 21

# Selection / view
$ws.Range("A1").Select()
$wb.Windows.Item(1).TabSelected = $false
